$wb = $excel.ActiveWorkbook

$wsSteps = $wb.Worksheets.Item("TestSteps")
$wsData  = $wb.Worksheets.Item("TestData")

# Add the new test-data row (copy formatting from the row above, then set values)
$wsData.Range("A4:E4").Copy()
$wsData.Range("A5:E5").PasteSpecial(-4122)
$wsData.Range("A5").Value = 4
$wsData.Range("B5").Value = "Enter Text to delete Comptroller Source/Group Profile"
$wsData.Range("C5").Value = 3243

# Update selection on TestSteps (no longer the active/tab-selected sheet)
$wsSteps.Range("A12").Select()

# Make TestData the active sheet/tab and set its new selection
$wsData.Activate()
$wsData.Range("C5").Select()
